$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.5607738878
$ws.Range("C2").Value = -224.72799182
$ws.Range("D2").Value = -225.28876571

$ws.Range("B3").Value = -0.5693531088
$ws.Range("C3").Value = -224.66963371
$ws.Range("D3").Value = -225.23898681

$ws.Range("B4").Value = -0.5696708344
$ws.Range("C4").Value = -224.63670635
$ws.Range("D4").Value = -225.20637719
